# Update the "Number of lines/tx's" (column E) values on the
# "Bus connections" sheet for the first two branches (rows 2 and 3)
# from 7500 to 10000. The dependent reactance formulas in column C
# recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bus connections")
$ws.Range("E2").Value = 10000
$ws.Range("E3").Value = 10000

# Make "Bus connections" the active/selected sheet (was "Bus index"
# before), and narrow its selection down to the single cell E4.
$ws.Activate()
$ws.Range("E4").Select()
